$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 3-21 (columns D, L, M, N, O, P, S).
# All other columns (A,B,C,E,F,G,H,I,J,K,Q,R,T) are identical across these
# rows, so only these seven columns need to be rewritten per the diff.
$rows = @{
    3  = @{ D = 44532; L = "Primera"; M = 100; N = 10000; O = 10000; P = 10000; S = 5000 }
    4  = @{ D = 44532; L = "Segunda"; M = 100; N = 8000;  O = 8000;  P = 8000;  S = 4000 }
    5  = @{ D = 44609; L = "Primera"; M = 100; N = 6500;  O = 7000;  P = 6750;  S = 3375 }
    6  = @{ D = 44609; L = "Segunda"; M = 50;  N = 6000;  O = 6000;  P = 6000;  S = 3000 }
    7  = @{ D = 44559; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    8  = @{ D = 44559; L = "Segunda"; M = 100; N = 5000;  O = 5000;  P = 5000;  S = 2500 }
    9  = @{ D = 44602; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    10 = @{ D = 44602; L = "Segunda"; M = 100; N = 5000;  O = 5000;  P = 5000;  S = 2500 }
    11 = @{ D = 44910; L = "Primera"; M = 200; N = 7500;  O = 8000;  P = 7750;  S = 3875 }
    12 = @{ D = 44988; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    13 = @{ D = 45014; L = "Primera"; M = 100; N = 7000;  O = 7500;  P = 7250;  S = 3625 }
    14 = @{ D = 44617; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    15 = @{ D = 44195; L = "Primera"; M = 200; N = 3000;  O = 3500;  P = 3250;  S = 1625 }
    16 = @{ D = 44195; L = "Segunda"; M = 100; N = 2500;  O = 2500;  P = 2500;  S = 1250 }
    17 = @{ D = 44944; L = "Primera"; M = 100; N = 7000;  O = 8000;  P = 7500;  S = 3750 }
    18 = @{ D = 44574; L = "Primera"; M = 200; N = 7000;  O = 8000;  P = 7500;  S = 3750 }
    19 = @{ D = 44574; L = "Segunda"; M = 100; N = 6000;  O = 6000;  P = 6000;  S = 3000 }
    20 = @{ D = 44216; L = "Primera"; M = 200; N = 3500;  O = 4000;  P = 3750;  S = 1875 }
    21 = @{ D = 44216; L = "Segunda"; M = 100; N = 3000;  O = 3000;  P = 3000;  S = 1500 }
}

foreach ($r in $rows.Keys) {
    $v = $rows[$r]
    $ws.Range("D$r").Value = $v.D
    $ws.Range("L$r").Value = $v.L
    $ws.Range("M$r").Value = $v.M
    $ws.Range("N$r").Value = $v.N
    $ws.Range("O$r").Value = $v.O
    $ws.Range("P$r").Value = $v.P
    $ws.Range("S$r").Value = $v.S
}
